{"js": "const replacements = [\n  { oldText: \"18\u00d798=\", newText: \"95\u00d770=\" },\n  { oldText: \"81\u00d769=\", newText: \"13\u00d762=\" },\n  { oldText: \"87\u00d745=\", newText: \"24\u00d793=\" },\n  { oldText: \"89\u00d762=\", newText: \"27\u00d790=\" },\n  { oldText: \"26\u00d712=\", newText: \"78\u00d711=\" },\n  { oldText: \"12\u00d748=\", newText: \"42\u00d736=\" },\n  { oldText: \"42\u00d743=\", newText: \"73\u00d721=\" },\n  { oldText: \"53\u00d795=\", newText: \"47\u00d744=\" },\n  { oldText: \"93\u00d770=\", newText: \"98\u00d728=\" },\n  { oldText: \"39\u00d792=\", newText: \"29\u00d712=\" },\n  { oldText: \"72\u00d752=\", newText: \"28\u00d734=\" },\n  { oldText: \"27\u00d754=\", newText: \"54\u00d785=\" },\n  { oldText: \"42\u00d733=\", newText: \"22\u00d726=\" },\n  { oldText: \"59\u00d726=\", newText: \"62\u00d719=\" },\n  { oldText: \"88\u00d781=\", newText: \"73\u00d740=\" },\n  { oldText: \"92\u00d714=\", newText: \"63\u00d796=\" },\n  { oldText: \"14\u00d791=\", newText: \"72\u00d712=\" },\n  { oldText: \"23\u00d729=\", newText: \"57\u00d793=\" },\n  { oldText: \"21\u00d772=\", newText: \"79\u00d752=\" },\n  { oldText: \"25\u00d756=\", newText: \"74\u00d721=\" },\n  { oldText: \"11\u00d754=\", newText: \"59\u00d788=\" },\n  { oldText: \"26\u00d758=\", newText: \"76\u00d742=\" },\n  { oldText: \"22\u00d772=\", newText: \"25\u00d773=\" },\n  { oldText: \"76\u00d746=\", newText: \"54\u00d720=\" },\n  { oldText: \"84\u00d797=\", newText: \"65\u00d768=\" },\n];\n\nconst body = context.document.body;\n\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    for (let i = 0; i < results.items.length; i++) {\n      results.items[i].insertText(newText, Word.InsertLocation.replace);\n    }\n    await context.sync();\n  }\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ old = \"18\u00d798=\"; new = \"95\u00d770=\" },\n    @{ old = \"81\u00d769=\"; new = \"13\u00d762=\" },\n    @{ old = \"87\u00d745=\"; new = \"24\u00d793=\" },\n    @{ old = \"89\u00d762=\"; new = \"27\u00d790=\" },\n    @{ old = \"26\u00d712=\"; new = \"78\u00d711=\" },\n    @{ old = \"12\u00d748=\"; new = \"42\u00d736=\" },\n    @{ old = \"42\u00d743=\"; new = \"73\u00d721=\" },\n    @{ old = \"53\u00d795=\"; new = \"47\u00d744=\" },\n    @{ old = \"93\u00d770=\"; new = \"98\u00d728=\" },\n    @{ old = \"39\u00d792=\"; new = \"29\u00d712=\" },\n    @{ old = \"72\u00d752=\"; new = \"28\u00d734=\" },\n    @{ old = \"27\u00d754=\"; new = \"54\u00d785=\" },\n    @{ old = \"42\u00d733=\"; new = \"22\u00d726=\" },\n    @{ old = \"59\u00d726=\"; new = \"62\u00d719=\" },\n    @{ old = \"88\u00d781=\"; new = \"73\u00d740=\" },\n    @{ old = \"92\u00d714=\"; new = \"63\u00d796=\" },\n    @{ old = \"14\u00d791=\"; new = \"72\u00d712=\" },\n    @{ old = \"23\u00d729=\"; new = \"57\u00d793=\" },\n    @{ old = \"21\u00d772=\"; new = \"79\u00d752=\" },\n    @{ old = \"25\u00d756=\"; new = \"74\u00d721=\" },\n    @{ old = \"11\u00d754=\"; new = \"59\u00d788=\" },\n    @{ old = \"26\u00d758=\"; new = \"76\u00d742=\" },\n    @{ old = \"22\u00d772=\"; new = \"25\u00d773=\" },\n    @{ old = \"76\u00d746=\"; new = \"54\u00d720=\" },\n    @{ old = \"84\u00d797=\"; new = \"65\u00d768=\" },\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.old\n    $find.Replacement.Text = $pair.new\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($pair.old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2) | Out-Null\n}"}
